$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order row: D-100
$ws.Range("A2").Value = "D-100"
$ws.Range("B2").Value = "NTU"
$ws.Range("C2").Value = "chicken nugget, burger"
$ws.Range("D2").Value = "burger : spicy"
# "false" must land as literal text (not a Boolean) - force text entry with a
# leading quote (as a user would in Excel), then reset the cell style back to
# Normal so no extra formatting is left behind on the cell.
$ws.Range("E2").Value = "'false"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Value = "NEW"
$ws.Range("G2").Value = "Cash"

# New order row: D-101
$ws.Range("A3").Value = "D-101"
$ws.Range("B3").Value = "NTU"
$ws.Range("C3").Value = "Chicken tenders"
$ws.Range("D3").Value = "Chicken tenders : spicy"
$ws.Range("E3").Value = "'false"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Value = "NEW"
$ws.Range("G3").Value = "Cash"
